# "update logs and weekly ptt"
# Fill in newly logged work entries on the "Week 3" worksheet (rows 54, 60-67,
# 71-72), matching the DESCRIPTION / Activity Type / HRs columns, then leave
# the selection where the author last left off on each sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Week 1")
$ws2 = $wb.Worksheets.Item("Week 2")
$ws3 = $wb.Worksheets.Item("Week 3")

# --- Week 3: new log entries -------------------------------------------------

$ws3.Range("A54").Value = "Helping UDL with data streaming parsing"
$ws3.Range("B54").Value = "Streaming Parsing support for UDL"
$ws3.Range("C54").Value = 1.5

$ws3.Range("A60").Value = "Update PowerPoint for meeting with TechSafety BC"
$ws3.Range("B60").Value = "Project Documents"
$ws3.Range("C60").Value = 1

$ws3.Range("B61").Value = "Anomaly Detection Model"
$ws3.Range("A61").Value = "Anomaly detection model discussions"
$ws3.Range("C61").Value = 1

$ws3.Range("A62").Value = "Meeting with TechSafetyBC"
$ws3.Range("B62").Value = "Client Meetings"
$ws3.Range("C62").Value = 1

$ws3.Range("A63").Value = "Post meeting discussion with UDL and discussionon next steps"
$ws3.Range("B63").Value = "Client Meetings / Data Access"
$ws3.Range("C63").Value = 1.25

$ws3.Range("A64").Value = "Go over with Ryan how to start downloading data from SkySpark manually"
$ws3.Range("B64").Value = "Data Access"
$ws3.Range("C64").Value = 0.25

$ws3.Range("A65").Value = "Meeting with Scott to discuss proposed method"
$ws3.Range("B65").Value = "UBC Meeting"
$ws3.Range("C65").Value = 0.75

$ws3.Range("A66").Value = "Post meeting discussion"
$ws3.Range("B66").Value = "Internal Meeting"
$ws3.Range("C66").Value = 0.25

$ws3.Range("A67").Value = "Figuring out a way to easily view and label anomalies"
$ws3.Range("B67").Value = "Anomaly Labelling"
$ws3.Range("C67").Value = 3.5

$ws3.Range("A71").Value = "Built a Shiny app to support labelling anomalies"
$ws3.Range("B71").Value = "Anomaly Labelling"
$ws3.Range("C71").Value = 8.5

$ws3.Range("A72").Value = "Update Week 4 Status Presentation"
$ws3.Range("B72").Value = "Project Documents"
$ws3.Range("C72").Value = 0.5

# --- Restore each sheet's last-known selection/scroll position --------------

$ws1.Activate()
$ws1.Range("A49").Select()

$ws2.Activate()
$ws2.Range("A55").Select()

$ws3.Activate()
$ws3.Range("A73").Select()
